$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Paragraph "React-router..." gets 12pt (240 twips) space-before.
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("React-router", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Paragraphs.Item(1).SpaceBefore = 12

# ------------------------------------------------------------------------
# 2) Insert an (empty) _GoBack bookmark right after "...名をつけられたりする
#    ので便利。" (marks the last editing position, as Word does automatically).
#    A collapsed range sitting exactly at that paragraph's end (when the
#    following paragraph is empty) cannot be fed to Bookmarks.Add directly,
#    so we anchor on a temporary marker run, create the bookmark there, and
#    remove the marker again.
# ------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("名をつけられたりするので便利。", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Collapse(0)
$rng2.InsertAfter("__TMP_GOBACK__")

$markerRng = $d.Content
$markerRng.Find.Execute("__TMP_GOBACK__", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$markerRng.Collapse(1)
$d.Bookmarks.Add("_GoBack", $markerRng)

$cleanupRng = $d.Content
$cleanupRng.Find.Execute("__TMP_GOBACK__", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$cleanupRng.Delete()

# ------------------------------------------------------------------------
# 3) "Test suites..." paragraph loses the eastAsia-hint font applied to its
#    paragraph mark (pilcrow) formatting.
# ------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("Test suites", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $rng3.Paragraphs.Item(1)
$markRng = $para.Range
$markRng.Collapse(0)
$markRng.Font.Name = ""
